# Apply crypto price/volume updates for Sun Feb 26 17:21:05 UTC 2023 GitHub Actions run
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "23.209.65"
$ws.Range("E2").Value = "  +0.71%  "
$ws.Range("D3").Value = "1.602.49"
$ws.Range("E3").Value = "  +0.24%  "
$ws.Range("D4").NumberFormat = "@"  # keep as text, not a number
$ws.Range("D4").Value = "1.002"
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("E5").Value = "  +0.06%  "
$ws.Range("D6").NumberFormat = "@"  # keep as text, not a number
$ws.Range("D6").Value = "304.77"
$ws.Range("E6").Value = "  +0.85%  "
$ws.Range("D7").NumberFormat = "@"  # keep as text, not a number
$ws.Range("D7").Value = "0.3763"
$ws.Range("E7").Value = "  -0.39%  "
$ws.Range("D8").NumberFormat = "@"  # keep as text, not a number
$ws.Range("D8").Value = "52.78"
$ws.Range("E8").Value = "  +4.17%  "
$ws.Range("D9").NumberFormat = "@"  # keep as text, not a number
$ws.Range("D9").Value = "0.3610"
$ws.Range("E9").Value = "  -0.85%  "
$ws.Range("D10").NumberFormat = "@"  # keep as text, not a number
$ws.Range("D10").Value = "1.263"
$ws.Range("E10").Value = "  +0.95%  "
$ws.Range("D11").NumberFormat = "@"  # keep as text, not a number
$ws.Range("D11").Value = "1.001"
$ws.Range("E11").Value = "  +0.02%  "
$ws.Range("D12").NumberFormat = "@"  # keep as text, not a number
$ws.Range("D12").Value = "0.08131"
$ws.Range("E12").Value = "  -0.10%  "
$ws.Range("D13").NumberFormat = "@"  # keep as text, not a number
$ws.Range("D13").Value = "22.83"
$ws.Range("E13").Value = "  +2.19%  "
$ws.Range("D14").NumberFormat = "@"  # keep as text, not a number
$ws.Range("D14").Value = "6.588"
$ws.Range("E14").Value = "  +0.26%  "
$ws.Range("D15").NumberFormat = "@"  # keep as text, not a number
$ws.Range("D15").Value = "7.347"
$ws.Range("E15").Value = "  -0.26%  "
$ws.Range("D16").NumberFormat = "@"  # keep as text, not a number
$ws.Range("D16").Value = "0.00001243"
$ws.Range("E16").Value = "  +0.09%  "
$ws.Range("D17").Value = "1.602.14"
$ws.Range("E17").Value = "  +0.15%  "
$ws.Range("D18").NumberFormat = "@"  # keep as text, not a number
$ws.Range("D18").Value = "94.01"
$ws.Range("E18").Value = "  +2.25%  "
$ws.Range("D19").NumberFormat = "@"  # keep as text, not a number
$ws.Range("D19").Value = "0.06929"
$ws.Range("E19").Value = "  +1.30%  "
$ws.Range("E20").Value = "  -0.12%  "
$ws.Range("D21").NumberFormat = "@"  # keep as text, not a number
$ws.Range("D21").Value = "6.525"
$ws.Range("E21").Value = "  +0.24%  "
$ws.Range("D22").NumberFormat = "@"  # keep as text, not a number
$ws.Range("D22").Value = "1.005"
$ws.Range("E22").Value = "  +0.31%  "
$ws.Range("D23").NumberFormat = "@"  # keep as text, not a number
$ws.Range("D23").Value = "12.87"
$ws.Range("E23").Value = "  -1.06%  "
$ws.Range("D24").Value = "23.221.63"
$ws.Range("E24").Value = "  +0.74%  "
$ws.Range("B25").Value = "LidoDAOToken"
$ws.Range("C25").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D25").NumberFormat = "@"  # keep as text, not a number
$ws.Range("D25").Value = "3.055"
$ws.Range("E25").Value = "  +10.54%  "
$ws.Range("B26").Value = "Toncoin"
$ws.Range("C26").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D26").NumberFormat = "@"  # keep as text, not a number
$ws.Range("D26").Value = "2.424"
$ws.Range("E26").Value = "  +2.51%  "
$ws.Range("E27").Value = "  +0.38%  "
$ws.Range("D28").NumberFormat = "@"  # keep as text, not a number
$ws.Range("D28").Value = "150.59"
$ws.Range("E28").Value = "  +0.95%  "
$ws.Range("D29").NumberFormat = "@"  # keep as text, not a number
$ws.Range("D29").Value = "5.264"
$ws.Range("E29").Value = "  +0.32%  "
$ws.Range("D30").NumberFormat = "@"  # keep as text, not a number
$ws.Range("D30").Value = "134.99"
$ws.Range("E30").Value = "  +0.45%  "
$ws.Range("D31").NumberFormat = "@"  # keep as text, not a number
$ws.Range("D31").Value = "2.409"
$ws.Range("E31").Value = "  +2.22%  "
$ws.Range("D32").NumberFormat = "@"  # keep as text, not a number
$ws.Range("D32").Value = "6.743"
$ws.Range("E32").Value = "  -0.68%  "
$ws.Range("D33").Value = "1.782.82"
$ws.Range("E33").Value = "  +0.55%  "
$ws.Range("D34").NumberFormat = "@"  # keep as text, not a number
$ws.Range("D34").Value = "0.9521"
$ws.Range("E34").Value = "  -0.65%  "
$ws.Range("D35").NumberFormat = "@"  # keep as text, not a number
$ws.Range("D35").Value = "0.02766"
$ws.Range("E35").Value = "  +2.30%  "
$ws.Range("D36").NumberFormat = "@"  # keep as text, not a number
$ws.Range("D36").Value = "0.07421"
$ws.Range("E36").Value = "  -1.47%  "
$ws.Range("D37").NumberFormat = "@"  # keep as text, not a number
$ws.Range("D37").Value = "10.27"
$ws.Range("E37").Value = "  +1.15%  "
$ws.Range("D38").NumberFormat = "@"  # keep as text, not a number
$ws.Range("D38").Value = "0.2513"
$ws.Range("E38").Value = "  -0.23%  "
$ws.Range("D39").NumberFormat = "@"  # keep as text, not a number
$ws.Range("D39").Value = "6.099"
$ws.Range("E39").Value = "  -1.64%  "
$ws.Range("D40").NumberFormat = "@"  # keep as text, not a number
$ws.Range("D40").Value = "0.08746"
$ws.Range("E40").Value = "  -0.80%  "
$ws.Range("D41").NumberFormat = "@"  # keep as text, not a number
$ws.Range("D41").Value = "1.405"
$ws.Range("E41").Value = "  +3.51%  "
$ws.Range("D42").NumberFormat = "@"  # keep as text, not a number
$ws.Range("D42").Value = "0.7088"
$ws.Range("E42").Value = "  +0.84%  "
$ws.Range("D43").NumberFormat = "@"  # keep as text, not a number
$ws.Range("D43").Value = "12.42"
$ws.Range("E43").Value = "  +1.36%  "
$ws.Range("D44").NumberFormat = "@"  # keep as text, not a number
$ws.Range("D44").Value = "15.82"
$ws.Range("E44").Value = "  +3.78%  "
$ws.Range("D45").NumberFormat = "@"  # keep as text, not a number
$ws.Range("D45").Value = "0.6516"
$ws.Range("E45").Value = "  -0.98%  "
$ws.Range("D46").NumberFormat = "@"  # keep as text, not a number
$ws.Range("D46").Value = "2.322"
$ws.Range("E46").Value = "  +2.14%  "
$ws.Range("D47").NumberFormat = "@"  # keep as text, not a number
$ws.Range("D47").Value = "1.000"
$ws.Range("E47").Value = "  +0.05%  "
$ws.Range("D48").NumberFormat = "@"  # keep as text, not a number
$ws.Range("D48").Value = "4.011"
$ws.Range("E48").Value = "  +0.31%  "
$ws.Range("E49").Value = "  +1.56%  "
$ws.Range("E50").Value = "  +0.21%  "
$ws.Range("D51").NumberFormat = "@"  # keep as text, not a number
$ws.Range("D51").Value = "1.196"
$ws.Range("E51").Value = "  -2.04%  "
